$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Cases")

$base = "G:\E_Git_Traingings_Ecplise_workspace\E_POC2\Screenshots\"

# Row 2 (CP-27 block) - new screenshot for K2
$ws.Range("K2").Value = $base + "monthPerformance_19_06_2023_14_35_40.png"

# Row 3 (CP-29 defect id + screenshot)
$ws.Range("J3").Value = "CP-29"
$ws.Range("K3").Value = $base + "lowAttendancePerformance_19_06_2023_14_36_36.png"

# Row 4 - Grid performance screenshot
$ws.Range("K4").Value = $base + "GirdPerformance_19_06_2023_14_37_33.png"

# Row 5 - characteristic performance screenshot
$ws.Range("K5").Value = $base + "characteristicPerformance_19_06_2023_14_37_37.png"

# Row 6 - Marks obtained screenshot
$ws.Range("K6").Value = $base + "Marks_Obtained_19_06_2023_14_37_41.png"

# Row 7 - total marks screenshot
$ws.Range("K7").Value = $base + "total_Marks_19_06_2023_14_37_45.png"

# Row 8 - marks percentage screenshot
$ws.Range("K8").Value = $base + "marks_Percentage_19_06_2023_14_37_49.png"

# Row 9 - attendance percentage screenshot
$ws.Range("K9").Value = $base + "attendance_Percentage_19_06_2023_14_37_53.png"

# Row 10 - group marking screenshot
$ws.Range("K10").Value = $base + "Group_Marking_19_06_2023_14_37_56.png"
